$d = $word.ActiveDocument

# Update the date paragraph (no trailing CR -> keeps single paragraph)
$d.Paragraphs.Item(1).Range.Text = "2025-02-01 Saturday"

# Update the table cells (20 rows x 5 columns), in document order
$t = $d.Tables.Item(1)
$values = @(
    "94-52=",
    "70-17=",
    "67-35=",
    "23+46=",
    "67-63=",
    "59+14=",
    "48-7=",
    "71-44=",
    "71-33=",
    "37-4=",
    "33-18=",
    "97-50=",
    "80-3=",
    "30+65=",
    "5+9=",
    "28+66=",
    "65+1=",
    "54-11=",
    "74-68=",
    "65-6=",
    "91-3=",
    "96-51=",
    "98-88=",
    "99-58=",
    "27+30=",
    "7+38=",
    "14+28=",
    "60+17=",
    "97-66=",
    "54-36=",
    "86-4=",
    "45+6=",
    "28-20=",
    "90-49=",
    "78-33=",
    "15+8=",
    "68-64=",
    "32+1=",
    "66-14=",
    "56-46=",
    "32+45=",
    "68-35=",
    "70-45=",
    "93-56=",
    "54-43=",
    "93-62=",
    "79-61=",
    "6+65=",
    "55+18=",
    "50+32=",
    "56+3=",
    "25-15=",
    "52-42=",
    "74-54=",
    "41+22=",
    "26+35=",
    "27+14=",
    "93-12=",
    "97-34=",
    "79-27=",
    "99-28=",
    "53-46=",
    "19+55=",
    "75-41=",
    "68-18=",
    "43-39=",
    "7+52=",
    "87-6=",
    "57+34=",
    "4+85=",
    "51-11=",
    "81-48=",
    "43+32=",
    "61-13=",
    "50-9=",
    "69-13=",
    "67+7=",
    "2+21=",
    "31+31=",
    "88+11=",
    "23-3=",
    "66-44=",
    "24-7=",
    "21-21=",
    "54+9=",
    "85-64=",
    "90-50=",
    "88-18=",
    "45-14=",
    "95-22=",
    "55+28=",
    "98-73=",
    "96-11=",
    "75-53=",
    "1+7=",
    "69-13=",
    "19+19=",
    "21+59=",
    "13+22=",
    "17+39="
)

$idx = 0
for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}

Write-Host "Updated $idx cells"
